# Updated cryptos list on Wed Oct 25 20:58:29 UTC 2023 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain text (e.g. "34.731.23", "0.555") that
# would otherwise be auto-coerced to a Number by the COM Value setter
# whenever it looks like a simple decimal. Force text interpretation for
# the data range first, write all values, then restore the default style
# so no stray per-cell formatting is left behind.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2 - Bitcoin
$ws.Range("D2").Value = "34.697.08"
$ws.Range("E2").Value = "  +2.92%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.788.40"
$ws.Range("E3").Value = "  +0.81%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.11%  "

# Row 5 - BNB
$ws.Range("D5").Value = "222.72"
$ws.Range("E5").Value = "  -0.98%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.554"
$ws.Range("E6").Value = "  -0.69%  "

# Row 7 - USDC
$ws.Range("D7").Value = "1.00"

# Row 8 - Solana
$ws.Range("D8").Value = "32.55"
$ws.Range("E8").Value = "  +7.68%  "

# Row 9 - Cardano
$ws.Range("D9").Value = "0.279"
$ws.Range("E9").Value = "  +0.73%  "

# Row 10 - Dogecoin
$ws.Range("D10").Value = "0.0683"
$ws.Range("E10").Value = "  +3.09%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.46%  "

# Row 12 - WrappedliquidstakedEther2.0
$ws.Range("D12").Value = "2.047.86"
$ws.Range("E12").Value = "  +0.80%  "

# Row 13 - Chainlink
$ws.Range("D13").Value = "11.01"
$ws.Range("E13").Value = "  +10.25%  "

# Row 14 - WrappedEther
$ws.Range("D14").Value = "1.792.75"
$ws.Range("E14").Value = "  +0.72%  "

# Row 15 - WrappedBTC
$ws.Range("D15").Value = "34.711.64"
$ws.Range("E15").Value = "  +2.96%  "

# Row 16 - Polygon
$ws.Range("D16").Value = "0.629"
$ws.Range("E16").Value = "  +1.15%  "

# Row 17 - Polkadot
$ws.Range("D17").Value = "4.30"
$ws.Range("E17").Value = "  +3.15%  "

# Row 18 - Litecoin
$ws.Range("D18").Value = "68.37"
$ws.Range("E18").Value = "  -0.04%  "

# Row 19 - BitcoinCash
$ws.Range("D19").Value = "252.55"
$ws.Range("E19").Value = "  +0.95%  "

# Row 20 - ShibaInu
$ws.Range("D20").Value = "0.0₃0779"
$ws.Range("E20").Value = "  +5.94%  "

# Row 21 - Dai
$ws.Range("D21").Value = "0.998"
$ws.Range("E21").Value = "  -0.18%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "10.46"
$ws.Range("E22").Value = "  +2.28%  "

# Row 23 - Uniswap
$ws.Range("D23").Value = "4.17"
$ws.Range("E23").Value = "  +0.36%  "

# Row 24 - Toncoin
$ws.Range("E24").Value = "  -0.81%  "

# Row 25 - Monero
$ws.Range("D25").Value = "158.83"
$ws.Range("E25").Value = "  +0.21%  "

# Row 26 - EthereumClassic
$ws.Range("D26").Value = "16.34"
$ws.Range("E26").Value = "  -0.26%  "

# Row 27 - Cosmos
$ws.Range("D27").Value = "7.02"
$ws.Range("E27").Value = "  +1.26%  "

# Row 28 - Stellar
$ws.Range("E28").Value = "  +0.19%  "

# Row 29 - BinanceUSD
$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  -0.05%  "

# Row 30 - Hedera
$ws.Range("D30").Value = "0.0514"
$ws.Range("E30").Value = "  +0.18%  "

# Row 31 - Filecoin
$ws.Range("D31").Value = "3.73"
$ws.Range("E31").Value = "  -1.43%  "

# Row 32 - PancakeSwap
$ws.Range("E32").Value = "  +0.00%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").Value = "3.56"
$ws.Range("E33").Value = "  +0.39%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  +2.15%  "

# Row 35 - Maker
$ws.Range("D35").Value = "1.433.19"
$ws.Range("E35").Value = "  -2.84%  "

# Row 36 - TrustWalletToken
$ws.Range("D36").Value = "1.05"
$ws.Range("E36").Value = "  -1.16%  "

# Row 37/38 - swap ImmutableX <-> VeChain (with refreshed price/volume)
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "0.0189"
$ws.Range("E37").Value = "  +2.39%  "

$ws.Range("B38").Value = "ImmutableX"
$ws.Range("C38").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D38").Value = "0.630"
$ws.Range("E38").Value = "  +0.55%  "

# Row 39 - Aave
$ws.Range("D39").Value = "82.62"
$ws.Range("E39").Value = "  -0.10%  "

# Row 40 - MXToken
$ws.Range("D40").Value = "2.80"
$ws.Range("E40").Value = "  +4.18%  "

# Row 41 - HuobiToken
$ws.Range("E41").Value = "  -0.07%  "

# Row 42 - ARBITRUM
$ws.Range("D42").Value = "0.900"
$ws.Range("E42").Value = "  +2.17%  "

# Row 43 - RenderToken
$ws.Range("E43").Value = "  -0.56%  "

# Row 44 - WEMIXToken
$ws.Range("E44").Value = "  -0.72%  "

# Row 45 - Kaspa
$ws.Range("E45").Value = "  -0.92%  "

# Row 46 - FraxShare
$ws.Range("D46").Value = "5.93"
$ws.Range("E46").Value = "  +4.39%  "

# Row 47 - RocketPoolETH
$ws.Range("D47").Value = "1.943.09"
$ws.Range("E47").Value = "  +0.73%  "

# Row 48 - Quant
$ws.Range("D48").Value = "104.35"
$ws.Range("E48").Value = "  +7.52%  "

# Row 49/50 - swap InjectiveProtocol <-> PaxDollar (with refreshed price/volume)
$ws.Range("B49").Value = "PaxDollar"
$ws.Range("C49").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D49").Value = "0.999"
$ws.Range("E49").Value = "  -0.11%  "

$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").Value = "11.93"
$ws.Range("E50").Value = "  +1.09%  "

# Row 51 - BitcoinSV
$ws.Range("D51").Value = "49.54"
$ws.Range("E51").Value = "  -2.48%  "

# Restore the default (unstyled) look for the Price column now that all
# values are safely stored as text.
$ws.Range("D2:D51").Style = "Normal"
